$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Link" column header
$ws.Range("I1").Value = "Link"

# New links in column I, written in the same order the original author
# appears to have entered them (this also controls shared-string order).
$ws.Range("I30").Value = "https://www.conrad.de/de/tracopower-tme-2405s-dcdc-wandler-print-24-vdc-5-vdc-200-ma-1-w-anzahl-ausgaenge-1-x-154477.html"
$ws.Range("I37").Value = "https://www.conrad.de/de/finder-361190054011-printrelais-5-vdc-10-a-1-wechsler-1-st-503243.html"
$ws.Range("I33").Value = "https://www.conrad.de/de/dip-schalter-polzahl-8-smd-apem-ikh0803000-1-st-700772.html"
$ws.Range("I25").Value = "https://www.conrad.de/de/tvs-diode-nexperia-pesd2can215-sot-23-262-v-230-w-1096200.html"
$ws.Range("I23").Value = "https://www.conrad.de/de/tvs-diode-bourns-smaj30ca-do-214ac-333-v-400-w-1056305.html"
$ws.Range("I21").Value = "https://www.conrad.de/de/quarzkristall-euroquartz-quarz-tc26-zylinder-32768-khz-125-pf-o-x-h-2-mm-x-62-mm-1-st-156007.html"
$ws.Range("I20").Value = "https://de.rs-online.com/web/p/quarzmodule/1710659/"
$ws.Range("I22").Value = "https://www.conrad.de/de/panjit-schottky-diode-gleichrichter-sr36-do-214aa-60-v-einzeln-1304924.html"
$ws.Range("I24").Value = "//"

# Row 21 (data row 22 counting the header) now references both D3-1 and D3-3,
# with quantity combined from the two original rows.
$ws.Range("D22").Value = "D3-1, D3-3"
$ws.Range("E22").Value = 2

$ws.Range("I32").Value = "https://www.conrad.de/de/pmic-waermemanagement-maxim-integrated-max31865atp-extern-spi-tqfn-20-ep-5x5-1123421.html"
$ws.Range("I34").Value = "https://de.rs-online.com/web/p/metalloxid-varistoren/7606961/"
$ws.Range("I35").Value = "https://de.rs-online.com/web/p/rueckstellende-sicherungen-smd/6478342/"
$ws.Range("I29").Value = "https://www.conrad.de/de/tvs-diode-stmicroelectronics-dviulc6-4sc6y-sot-23-6l-6-v-1183886.html"
$ws.Range("I19").Value = "https://www.conrad.de/de/yageo-cc1206zpy5v7bb475-keramik-kondensator-smd-1206-47-f-16-v-20-1-st-445372.html"
$ws.Range("I18").Value = "https://www.conrad.de/de/tantal-kondensator-smd-22-f-63-v-20-l-x-b-35-mm-x-28-mm-panasonic-6tpu22msi-1-st-1479548.html"
$ws.Range("I17").Value = "https://de.rs-online.com/web/p/keramik-multilayer-kondensatoren/9159328/"
$ws.Range("I10").Value = "https://www.conrad.de/de/microtech-cmf0603402r0110-duennschicht-widerstand-402-smd-0603-01-w-01-10-ppm-1-st-1457041.html"

# Row 24 (D3-3) had its own "Bestückung" entry merged into row 22 above, and
# its quantity is now 0 (fully absorbed).
$ws.Range("D24").ClearContents()
$ws.Range("E24").Value = 0

# Update the view: scroll down and select F29, matching the author's
# final cursor position when they saved the workbook.
$ws.Range("F29").Select() | Out-Null

